$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29
$ws.Range("A3").Value = "Média"
$ws.Range("B3").Value = 29
